$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Oppvisningsstevne Iskanten 9 februar 2026'
$ws.Range("D4").Value = '17:04:00'
$ws.Range("B5").Value = '17:04:00'
$ws.Range("D5").Value = '17:07:44'
$ws.Range("E5").Value = 'Elena Sophia Sandnes-Strømmen'
$ws.Range("F5").Value = 'Fana Idrettslag'
$ws.Range("B6").Value = '17:07:44'
$ws.Range("D6").Value = '17:11:28'
$ws.Range("E6").Value = 'Eira Olava Bortne Ludvigsen'
$ws.Range("F6").Value = 'Bergen Kunstløpklubb'
$ws.Range("B7").Value = '17:11:28'
$ws.Range("D7").Value = '17:15:12'
$ws.Range("E7").Value = 'Sara Barbro Kyte'
$ws.Range("F7").Value = 'Bergen Kunstløpklubb'
$ws.Range("B8").Value = '17:15:12'
$ws.Range("D8").Value = '17:18:56'
$ws.Range("E8").Value = 'Mie Mariell Sævereid'
$ws.Range("B9").Value = '17:18:56'
$ws.Range("D9").Value = '17:22:40'
$ws.Range("E9").Value = 'Amanda Ansnes Lima'
$ws.Range("B10").Value = '17:22:40'
$ws.Range("D10").Value = '17:26:24'
$ws.Range("E10").Value = 'Emilie Morseth'
$ws.Range("B11").Value = '17:26:24'
$ws.Range("D11").Value = '17:30:08'
$ws.Range("E11").Value = 'Leah Kalvik'
$ws.Range("F11").Value = 'Loddefjord IL'
$ws.Range("B12").Value = '17:30:08'
$ws.Range("D12").Value = '17:33:52'
$ws.Range("E12").Value = 'Aurelia Landschulze'
$ws.Range("F12").Value = 'Fana Idrettslag'
$ws.Range("B13").Value = 'ca. 17:33:52'
$ws.Range("D13").Value = '17:37:52'
$ws.Range("B14").Value = '17:37:52'
$ws.Range("D14").Value = '17:41:36'
$ws.Range("E14").Value = 'Sarolt Szofia Papdi'
$ws.Range("F14").Value = 'Bergen Kunstløpklubb'
$ws.Range("B15").Value = '17:41:36'
$ws.Range("D15").Value = '17:45:20'
$ws.Range("E15").Value = 'Camilla Tveit'
$ws.Range("B16").Value = '17:45:20'
$ws.Range("D16").Value = '17:49:04'
$ws.Range("E16").Value = 'Hanna Wangsuk Tveita'
$ws.Range("F16").Value = 'Loddefjord IL'
$ws.Range("B17").Value = '17:49:04'
$ws.Range("D17").Value = '17:52:48'
$ws.Range("E17").Value = 'Patricija Levickaite'
$ws.Range("B18").Value = '17:52:48'
$ws.Range("D18").Value = '17:56:32'
$ws.Range("E18").Value = 'Hennie Markestad'
$ws.Range("F18").Value = 'Bergen Kunstløpklubb'
$ws.Range("B19").Value = '17:56:32'
$ws.Range("D19").Value = '18:00:16'
$ws.Range("B20").Value = '18:00:16'
$ws.Range("D20").Value = '18:04:00'
$ws.Range("E20").Value = 'Aylin Morseth'
$ws.Range("F20").Value = 'Fana Idrettslag'
$ws.Range("B21").Value = '18:04:00'
$ws.Range("D21").Value = '18:07:44'
$ws.Range("E21").Value = 'Anne Kristoffersen'
$ws.Range("B22").Value = 'ca. 18:07:44'
$ws.Range("D22").Value = '18:11:44'
$ws.Range("B23").Value = '18:11:44'
$ws.Range("D23").Value = '18:15:28'
$ws.Range("E23").Value = 'Angela Chen'
$ws.Range("F23").Value = 'Fana Idrettslag'
$ws.Range("B24").Value = '18:15:28'
$ws.Range("D24").Value = '18:19:12'
$ws.Range("E24").Value = 'Frida Lovisa Østerberg'
$ws.Range("F24").Value = 'Bergen Kunstløpklubb'
$ws.Range("B25").Value = '18:19:12'
$ws.Range("D25").Value = '18:22:56'
$ws.Range("E25").Value = 'Eleanora Egle'
$ws.Range("F25").Value = 'Loddefjord IL'
$ws.Range("B26").Value = '18:22:56'
$ws.Range("D26").Value = '18:26:40'
$ws.Range("E26").Value = 'Frida Qianlu He'
$ws.Range("F26").Value = 'Loddefjord IL'
$ws.Range("B27").Value = '18:26:40'
$ws.Range("D27").Value = '18:30:24'
$ws.Range("E27").Value = 'Frida Pasko Hansen'
$ws.Range("F27").Value = 'Loddefjord IL'
$ws.Range("B28").Value = '18:30:24'
$ws.Range("D28").Value = '18:34:08'
$ws.Range("E28").Value = 'Aksel Eriksen'
$ws.Range("B29").Value = '18:34:08'
$ws.Range("D29").Value = '18:37:52'
$ws.Range("E29").Value = 'Valentina Pinker-Spilde'
$ws.Range("F29").Value = 'Fana Idrettslag'
$ws.Range("B30").Value = '18:37:52'
$ws.Range("D30").Value = '18:41:36'
$ws.Range("E30").Value = 'Mille Isabell Steen Rein'
$ws.Range("F30").Value = 'Loddefjord IL'
$ws.Range("A32").Value = 'Generert 31.01.2026 23:53'
